# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Sat Aug 17 11:32:45 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.149.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.45%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.588.10'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.37%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.29%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.25'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.08%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.24%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.564'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.70%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.600.88'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.86%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.42'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.94%  '

# Row 11
$ws.Range("E11").Value = '  -0.13%  '

# Row 12
$ws.Range("E12").Value = '  -2.65%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.042.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.48%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.061.45'
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.49'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.53%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.592.98'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.52%  '

# Row 18
$ws.Range("E18").Value = '  -0.88%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '344.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.11%  '

# Row 20
$ws.Range("E20").Value = '  -0.38%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.54%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.40%  '

# Row 23
$ws.Range("E23").Value = '  +0.01%  '

# Row 24
$ws.Range("E24").Value = '  +2.72%  '

# Row 25
$ws.Range("E25").Value = '  -0.54%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.405'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.14%  '

# Row 27
$ws.Range("E27").Value = '  +0.21%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.07%  '

# Row 29
$ws.Range("E29").Value = '  +0.10%  '

# Row 30
$ws.Range("E30").Value = '  -3.48%  '

# Row 31
$ws.Range("E31").Value = '  +1.24%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.87'
$ws.Range("D32").Style = "Normal"

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.71'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.35%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.27%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.96'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.17%  '

# Row 36
$ws.Range("E36").Value = '  -1.50%  '

# Row 37
$ws.Range("E37").Value = '  +1.44%  '

# Row 38
$ws.Range("E38").Value = '  +0.98%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.824'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.07%  '

# Row 40
$ws.Range("E40").Value = '  -5.25%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.52'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.39%  '

# Row 42
$ws.Range("E42").Value = '  +0.27%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.599'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.57%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.79'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.01%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '268.34'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.09%  '

# Row 46
$ws.Range("E46").Value = '  -0.25%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0514'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.52%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.23%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.958.69'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.48%  '

# Row 50
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0221'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.71%  '

# Row 51
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.38%  '
